# Applies the commit "Horarios actualizados Linea 141 - 287" update to
# horarios-141-2026-01-13.xlsx: refreshes the scrape timestamp / row counts on
# each of the 3 sheets and writes the newly-scraped rows (the scraper re-sorts
# the table by arrival time each run, so a handful of older rows also shift).

$wb = $excel.ActiveWorkbook

function Set-Row($ws, $r, $horaScrap, $horaLlegada, $linea, $minutos, $parada) {
    $ws.Cells.Item($r, 1).Value = $horaScrap
    $ws.Cells.Item($r, 2).Value = $horaLlegada
    $ws.Cells.Item($r, 3).Value = $linea
    $ws.Cells.Item($r, 4).Value = $minutos
    $ws.Cells.Item($r, 5).Value = $parada
}

$lastUpdate = "17:59:03"

## --- Sheet 1: LP1912 (was A1:E421 / 416 rows -> A1:E446 / 441 rows) ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: $lastUpdate"
$ws1.Cells.Item(3,1).Value = "Total filas: 441"

Set-Row $ws1 253 '13:53:08' '13:57' '16_SANTA ANA' 4 'LP1912'
Set-Row $ws1 254 '12:37:14' '13:57' '16_P MOR-167 Y 521' 80 'LP1912'
Set-Row $ws1 281 '13:19:56' '15:04' '10_OLMOS' 105 'LP1912'
Set-Row $ws1 282 '14:46:52' '15:04' '23_HERNANDEZ' 18 'LP1912'
Set-Row $ws1 362 '17:14:55' '17:21' '10_OLMOS' 7 'LP1912'
Set-Row $ws1 364 '15:31:33' '17:21' '26_HERNANDEZ' 110 'LP1912'
Set-Row $ws1 392 '17:59:03' '17:59' '16_SANTA ANA' 0 'LP1912'
Set-Row $ws1 393 '17:41:19' '18:01' '16_SANTA ANA' 20 'LP1912'
Set-Row $ws1 394 '17:59:03' '18:02' '10_OLMOS' 3 'LP1912'
Set-Row $ws1 395 '16:33:08' '18:04' '17_ROMERO' 91 'LP1912'
Set-Row $ws1 396 '17:14:55' '18:04' '23_HERNANDEZ' 50 'LP1912'
Set-Row $ws1 397 '17:41:19' '18:05' '23_HERNANDEZ' 24 'LP1912'
Set-Row $ws1 398 '17:59:03' '18:06' '23_HERNANDEZ' 7 'LP1912'
Set-Row $ws1 399 '17:59:03' '18:06' '17_ROMERO' 7 'LP1912'
Set-Row $ws1 400 '17:14:55' '18:08' '14_ABASTO' 54 'LP1912'
Set-Row $ws1 401 '16:53:01' '18:09' '14_ABASTO' 76 'LP1912'
Set-Row $ws1 402 '17:41:19' '18:10' '14_ABASTO' 29 'LP1912'
Set-Row $ws1 403 '17:59:03' '18:12' '14_ABASTO' 13 'LP1912'
Set-Row $ws1 404 '17:59:03' '18:13' '16_SANTA ANA' 14 'LP1912'
Set-Row $ws1 405 '17:14:55' '18:16' '15_ABASTO' 62 'LP1912'
Set-Row $ws1 406 '17:41:19' '18:16' '10_OLMOS' 35 'LP1912'
Set-Row $ws1 407 '17:59:03' '18:18' '15_ABASTO' 19 'LP1912'
Set-Row $ws1 408 '17:59:03' '18:18' '10_OLMOS' 19 'LP1912'
Set-Row $ws1 409 '17:14:55' '18:20' '26_HERNANDEZ' 66 'LP1912'
Set-Row $ws1 410 '16:33:08' '18:21' '26_HERNANDEZ' 108 'LP1912'
Set-Row $ws1 411 '17:59:03' '18:23' '26_HERNANDEZ' 24 'LP1912'
Set-Row $ws1 412 '17:59:03' '18:25' '14_ABASTO' 44 'LP1912'
Set-Row $ws1 413 '17:59:03' '18:26' '14_ABASTO' 27 'LP1912'
Set-Row $ws1 414 '17:14:55' '18:27' '215C_EL PATO' 73 'LP1912'
Set-Row $ws1 415 '16:33:08' '18:28' '215C_EL PATO' 115 'LP1912'
Set-Row $ws1 416 '17:59:03' '18:30' '215C_EL PATO' 31 'LP1912'
Set-Row $ws1 417 '17:14:55' '18:31' '11X44_ETCHEVERRY' 77 'LP1912'
Set-Row $ws1 418 '16:43:37' '18:32' '11X44_ETCHEVERRY' 109 'LP1912'
Set-Row $ws1 419 '17:59:03' '18:34' '11X44_ETCHEVERRY' 35 'LP1912'
Set-Row $ws1 420 '17:41:19' '18:40' '15_ABASTO' 59 'LP1912'
Set-Row $ws1 421 '17:59:03' '18:42' '15_ABASTO' 43 'LP1912'
Set-Row $ws1 422 '17:14:55' '18:47' '14X44_ABASTO' 93 'LP1912'
Set-Row $ws1 423 '16:53:01' '18:48' '14X44_ABASTO' 115 'LP1912'
Set-Row $ws1 424 '17:59:03' '18:50' '14X44_ABASTO' 51 'LP1912'
Set-Row $ws1 425 '17:14:55' '18:58' '215A_EL PATO' 104 'LP1912'
Set-Row $ws1 426 '17:41:19' '18:59' '215A_EL PATO' 78 'LP1912'
Set-Row $ws1 427 '17:59:03' '19:00' '215A_EL PATO' 61 'LP1912'
Set-Row $ws1 428 '17:14:55' '19:04' '11_ETCHEVERRY' 110 'LP1912'
Set-Row $ws1 429 '17:41:19' '19:05' '11_ETCHEVERRY' 84 'LP1912'
Set-Row $ws1 430 '17:41:19' '19:05' '23_HERNANDEZ' 84 'LP1912'
Set-Row $ws1 431 '17:59:03' '19:06' '23_HERNANDEZ' 67 'LP1912'
Set-Row $ws1 432 '17:59:03' '19:06' '11_ETCHEVERRY' 67 'LP1912'
Set-Row $ws1 433 '17:14:55' '19:10' '16_P MOR-SANTA ANA' 116 'LP1912'
Set-Row $ws1 434 '17:41:19' '19:11' '16_P MOR-SANTA ANA' 90 'LP1912'
Set-Row $ws1 435 '17:59:03' '19:12' '16_P MOR-SANTA ANA' 73 'LP1912'
Set-Row $ws1 436 '17:41:19' '19:17' '27_EL RETIRO' 96 'LP1912'
Set-Row $ws1 437 '17:59:03' '19:18' '27_EL RETIRO' 79 'LP1912'
Set-Row $ws1 438 '17:41:19' '19:21' '26_HERNANDEZ' 100 'LP1912'
Set-Row $ws1 439 '17:59:03' '19:23' '26_HERNANDEZ' 84 'LP1912'
Set-Row $ws1 440 '17:41:19' '19:30' '225_GOMEZ' 109 'LP1912'
Set-Row $ws1 441 '17:59:03' '19:32' '225_GOMEZ' 93 'LP1912'
Set-Row $ws1 442 '17:41:19' '19:40' '215C_EL PATO' 119 'LP1912'
Set-Row $ws1 443 '17:59:03' '19:42' '215C_EL PATO' 103 'LP1912'
Set-Row $ws1 444 '17:59:03' '19:52' '16_P MOR-SANTA ANA' 113 'LP1912'
Set-Row $ws1 445 '17:59:03' '19:52' '11X44_ETCHEVERRY' 113 'LP1912'
Set-Row $ws1 446 '17:59:03' '19:53' '81_EL PELIGRO' 114 'LP1912'

## --- Sheet 2: LP1912-215 (was A1:E50 / 45 rows -> A1:E53 / 48 rows) ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: $lastUpdate"
$ws2.Cells.Item(3,1).Value = "Total filas: 48"

Set-Row $ws2 48 '17:59:03' '18:30' '215C_EL PATO' 31 'LP1912'
Set-Row $ws2 49 '17:14:55' '18:58' '215A_EL PATO' 104 'LP1912'
Set-Row $ws2 50 '17:41:19' '18:59' '215A_EL PATO' 78 'LP1912'
Set-Row $ws2 51 '17:59:03' '19:00' '215A_EL PATO' 61 'LP1912'
Set-Row $ws2 52 '17:41:19' '19:40' '215C_EL PATO' 119 'LP1912'
Set-Row $ws2 53 '17:59:03' '19:42' '215C_EL PATO' 103 'LP1912'

## --- Sheet 3: 6203-6173 (was A1:E62 / 57 rows -> A1:E63 / 58 rows) ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: $lastUpdate"
$ws3.Cells.Item(3,1).Value = "Total filas: 58"

Set-Row $ws3 63 '17:59:03' '19:54' '215C_LA PLATA' 115 'L6203'

